# Helper: pack an (r,g,b) triple into the VBA-style BGR-packed long used by
# the PowerPoint COM "RGB" values (R + G*256 + B*65536).
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 16: switch its table style to the new built-in style id.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(16)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{70D2011A-B664-471B-A808-F2624AC4DB5D}")
    }
}

# ---------------------------------------------------------------------------
# 2) Presentation theme: change the deck's colour scheme from the "Integral"
#    palette over to the classic "Office Theme" palette.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$officeThemeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i]
    $colorScheme.Item($i + 1).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
